$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking "Price" values stay as Text (matching original inlineStr cells)
# by forcing a Text number format before assignment, then resetting the style so no
# stray style index is left attached to the cell (matches the source which has no "s" attr).
$priceCells = @("D2", "D3", "D5", "D7", "D8", "D10", "D12", "D13", "D14", "D17", "D18", "D19", "D20", "D22", "D24", "D25", "D26", "D29", "D30", "D34", "D39", "D40", "D41", "D42", "D44", "D46", "D48", "D49", "D50", "D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '37.774.14'
$ws.Range('D3').Value = '2.077.78'
$ws.Range('D5').Value = '232.71'
$ws.Range('D7').Value = '1.00'
$ws.Range('D8').Value = '57.13'
$ws.Range('D10').Value = '0.0784'
$ws.Range('D12').Value = '2.373.30'
$ws.Range('D13').Value = '14.43'
$ws.Range('D14').Value = '20.93'
$ws.Range('D17').Value = '2.086.37'
$ws.Range('D18').Value = '37.695.83'
$ws.Range('D19').Value = '6.13'
$ws.Range('D20').Value = '70.55'
$ws.Range('D22').Value = '228.20'
$ws.Range('D24').Value = '2.40'
$ws.Range('D25').Value = '2.38'
$ws.Range('D26').Value = '170.95'
$ws.Range('D29').Value = '1.44'
$ws.Range('D30').Value = '19.36'
$ws.Range('D34').Value = '4.60'
$ws.Range('D39').Value = '5.41'
$ws.Range('D40').Value = '0.0999'
$ws.Range('D41').Value = '2.94'
$ws.Range('D42').Value = '97.85'
$ws.Range('D44').Value = '1.447.55'
$ws.Range('D46').Value = '4.18'
$ws.Range('D48').Value = '15.61'
$ws.Range('D49').Value = '7.42'
$ws.Range('D50').Value = '3.01'
$ws.Range('D51').Value = '2.269.82'

foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).Style = "Normal"
}

# Volume(1h) percentage cells are already non-numeric text (leading/trailing spaces, "%")
# so a direct Value assignment keeps them as Text.
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('E6').Value = '  +1.06%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +1.03%  '
$ws.Range('E9').Value = '  +1.40%  '
$ws.Range('E10').Value = '  +2.95%  '
$ws.Range('E11').Value = '  +3.06%  '
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('E14').Value = '  +1.20%  '
$ws.Range('E15').Value = '  -2.05%  '
$ws.Range('E16').Value = '  +2.93%  '
$ws.Range('E17').Value = '  +0.92%  '
$ws.Range('E18').Value = '  +1.30%  '
$ws.Range('E19').Value = '  -3.98%  '
$ws.Range('E20').Value = '  +1.59%  '
$ws.Range('E21').Value = '  +1.09%  '
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('E24').Value = '  -1.08%  '
$ws.Range('E25').Value = '  -0.48%  '
$ws.Range('E26').Value = '  +2.67%  '
$ws.Range('E27').Value = '  +10.18%  '
$ws.Range('E28').Value = '  +1.79%  '
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('E30').Value = '  +2.62%  '
$ws.Range('E31').Value = '  +1.04%  '
$ws.Range('E32').Value = '  +3.09%  '
$ws.Range('E33').Value = '  +1.59%  '
$ws.Range('E34').Value = '  +0.18%  '
$ws.Range('E35').Value = '  +0.86%  '
$ws.Range('E36').Value = '  +4.00%  '
$ws.Range('E37').Value = '  +4.97%  '
$ws.Range('E38').Value = '  -0.18%  '
$ws.Range('E39').Value = '  -4.00%  '
$ws.Range('E40').Value = '  +7.90%  '
$ws.Range('E41').Value = '  -0.40%  '
$ws.Range('E42').Value = '  +2.05%  '
$ws.Range('E43').Value = '  +1.18%  '
$ws.Range('E44').Value = '  -2.04%  '
$ws.Range('E45').Value = '  -0.28%  '
$ws.Range('E46').Value = '  -3.18%  '
$ws.Range('E47').Value = '  +3.01%  '
$ws.Range('E48').Value = '  +3.00%  '
$ws.Range('E49').Value = '  +3.76%  '
$ws.Range('E50').Value = '  +1.32%  '
$ws.Range('E51').Value = '  +0.48%  '
